$wb = $excel.ActiveWorkbook

# "Test Cases" overview sheet: just move the selection/active cell to C16
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("C16").Select()

# TestCase_A5: header D1 "PASS" -> "Results"
$ws = $wb.Worksheets.Item("TestCase_A5")
$ws.Range("D1").Value = "Results"

# TestCase_A6: header D1 "PASS" -> "Results"
$ws = $wb.Worksheets.Item("TestCase_A6")
$ws.Range("D1").Value = "Results"

# TestCase_A7: header D1 "PASS" -> "Results", selection moves D2 -> D1
$ws = $wb.Worksheets.Item("TestCase_A7")
$ws.Range("D1").Value = "Results"
$ws.Range("D1").Select()

# TestCase_A8: header D1 "PASS" -> "Results", selection moves D2 -> D1
$ws = $wb.Worksheets.Item("TestCase_A8")
$ws.Range("D1").Value = "Results"
$ws.Range("D1").Select()

# TestCase_A9: header D1 "PASS" -> "Results"
$ws = $wb.Worksheets.Item("TestCase_A9")
$ws.Range("D1").Value = "Results"

# TestCase_A10: header F1 "PASS" -> "Results", selection moves C15 -> F1
$ws = $wb.Worksheets.Item("TestCase_A10")
$ws.Range("F1").Value = "Results"
$ws.Range("F1").Select()

# TestCase_A11: header D1 "PASS" -> "Results"
$ws = $wb.Worksheets.Item("TestCase_A11")
$ws.Range("D1").Value = "Results"

# TestCase_A12: header F1 "PASS" -> "Results"
$ws = $wb.Worksheets.Item("TestCase_A12")
$ws.Range("F1").Value = "Results"

# TestCase_A19: header C1 "PASS" -> "Results", selection moves E3 -> C1
$ws = $wb.Worksheets.Item("TestCase_A19")
$ws.Range("C1").Value = "Results"
$ws.Range("C1").Select()

# Re-select the "Test Cases" sheet so it stays the tab-selected / active sheet
$wb.Worksheets.Item("Test Cases").Activate()
